$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")

# Insert a new blank column before column N (shifts old N:P -> O:Q)
$ws.Columns("N").Insert()

# Make "Repayment Schedule" the active sheet and select cell K19,
# which moves tabSelected from "Summary" to "Repayment Schedule".
$ws.Activate()
$ws.Range("K19").Select()
